$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H (col 8); old H shifts to I, old I shifts to J.
# This splits the old single "screenMotionCycle" field into Min (new col H)
# and Max (old col, now I) fields.
$ws.Columns(8).Insert()

# New column H => screenMotionCycleMin
$ws.Range("H1").Value = "screenMotionCycleMin"
$ws.Range("H2").Value = "모션 교체 주기 min (s)"
$ws.Range("H3").Value = "int"
$ws.Range("H4").Value = 5

# Shifted old column (now I) => screenMotionCycleMax (value 10 kept as-is)
$ws.Range("I1").Value = "screenMotionCycleMax"
$ws.Range("I2").Value = "모션 교체 주기 max (s)"

# Column widths: H=23.25, I=24.75 (already inherited), J=76.25 (chars equivalent
# accounting for the runtime's internal pixel rounding, MDW=7)
$ws.Columns(8).ColumnWidth = 22.571428571428573
$ws.Columns(10).ColumnWidth = 75.57142857142857

# Restore the cursor/selection position to match the edited workbook
$ws.Range("J10").Select()

Write-Host "done"
